$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","H","K","L","M","N","O")
$data = @(
  @(0.5971311412731666,0.1924353492520439,0.05269701357153878,0.1283788784727946,1.056645474702485,0.07973214163530429,0.2977978958269034,0.190895431617605,0.1672341601680039,2.058920891053223,3.807783358460938),
  @(0.5617636639494776,0.191872575847384,0.05072020775980945,0.1288173709433948,1.056773709657563,0.07973214163530429,0.2653995575083883,0.1883031288399621,0.1606412002114723,2.078227358966036,3.821891200383789),
  @(0.5402721185319876,0.1915311518114606,0.04949122407793993,0.129133568108367,1.057390381254947,0.07973214163530429,0.2455342934937903,0.1868017948968799,0.1566671286469159,2.090690007393814,3.832331535255662),
  @(0.5315711090177047,0.1913930744882997,0.04898659570405783,0.1292742507470468,1.057777067769464,0.07973214163530429,0.237446367137494,0.1862127748919349,0.155066396244159,2.095921657825531,3.837033474484059),
  @(0.5301297710824713,0.1913702111909501,0.04890257313589075,0.1292983260959915,1.057849457388812,0.07973214163530429,0.2361038284507799,0.1861163467264291,0.1548017308077547,2.096799612436002,3.837841259993866),
  @(0.5401545422301126,0.191529285351681,0.04948443387607426,0.1291354174775865,1.057395047905999,0.07973214163530429,0.2454251866440558,0.1867937588204356,0.1566454646001141,2.090759943645183,3.832393135265534),
  @(0.5848902359443571,0.1922404610967661,0.05201858002065762,0.1285203363893928,1.056578110902478,0.07973214163530429,0.2866214993039193,0.1899828919965358,0.1649456105347191,2.065451509692149,3.812278885906636),
  @(0.6743749815221918,0.1936669977813636,0.05686666414210606,0.1276859123111755,1.059240366617026,0.07973214163530429,0.3676102577206279,0.1969512646506502,0.1818053654811251,2.020646882275406,3.786931994741678),
  @(0.7411701001491906,0.1947336084718714,0.06035401061169665,0.1272984062128799,1.063792412862114,0.07973214163530429,0.4272222574185207,0.202504045389162,0.1945435030925751,1.990666286348572,3.776892552694846),
  @(0.7717809948949252,0.1952226733141984,0.06192420767510498,0.1271708716254949,1.066426347129209,0.07973214163530429,0.4543624251345761,0.2051237213645862,0.2004138683314949,1.977664267174015,3.77418700178751),
  @(0.7834044611736317,0.1954084061648587,0.06251645292167041,0.1271295672909645,1.067504647034895,0.07973214163530429,0.4646425558665044,0.2061291401725782,0.2026476106202253,1.97283221280332,3.773429925600226),
  @(0.7808997362413095,0.1953683818153351,0.06238900738140529,0.1271381523288042,1.067268820010881,0.07973214163530429,0.4624284316227829,0.2059120101628622,0.2021660576078546,1.973868810952705,3.773581083710098),
  @(0.7727366310930108,0.195237943067859,0.0619729793623236,0.1271673335071526,1.066513439003174,0.07973214163530429,0.4552081255426117,0.2052061695210767,0.2005974246849505,1.977264896621922,3.774119357669917),
  @(0.7677406139069376,0.1951581145924237,0.06171784291202442,0.1271861175924691,1.066061277133457,0.07973214163530429,0.4507858219684522,0.2047755653740921,0.199637989548819,1.979357017700188,3.774483890428684),
  @(0.7391740886672267,0.1947017228433339,0.06025106632875321,0.12730772011939,1.063631603329767,0.07973214163530429,0.4254489928840997,0.2023347228506793,0.1941613738473151,1.991528805311075,3.777106805330135),
  @(0.7217067441538063,0.1944227150798099,0.05934707739861977,0.1273947902348134,1.062285252098363,0.07973214163530429,0.4099110581730656,0.2008612925732649,0.1908209519803492,1.999158794896989,3.77919248503278),
  @(0.7116812490530435,0.1942626015551312,0.05882560206813991,0.1274494598057476,1.061563882147418,0.07973214163530429,0.4009761833241896,0.2000226368922569,0.1889067646117439,2.003607253207342,3.780567333431748),
  @(0.7082904593768831,0.1942084530649595,0.05864877828440029,0.1274687587541123,1.061328747183033,0.07973214163530429,0.3979513693472825,0.1997401998168584,0.1882598829400308,2.005123710121858,3.781062934362069),
  @(0.7235639771262754,0.1944523783688936,0.0594434665392356,0.1273850466515398,1.062423087278859,0.07973214163530429,0.4115648804272212,0.2010172291430052,0.1911758080665464,1.998340370879689,3.778952327753018),
  @(0.775133475396018,0.1952762417421496,0.0620952409528357,0.1271585727262341,1.066733118282968,0.07973214163530429,0.4573288339638282,0.2054131287788294,0.2010578791196664,1.976264899653259,3.77395399656848),
  @(0.8090220938199479,0.1958177908856698,0.06381459554159363,0.1270512945971021,1.070021392605767,0.07973214163530429,0.4872539609879993,0.2083642026055088,0.2075790563325768,1.962370813901517,3.772246143260304),
  @(0.7909183860223266,0.1955284782939373,0.06289820775998578,0.1271048299292552,1.068223276732951,0.07973214163530429,0.4712810758240096,0.2067820353483398,0.2040928886248423,1.969737522654251,3.773015093655914),
  @(0.7227242699898397,0.194438966686306,0.0593998944896299,0.1273894373591169,1.062360607937435,0.07973214163530429,0.4108171932975608,0.2009467039286932,0.1910153581039111,1.998710187638421,3.77906035534923),
  @(0.649980945936079,0.1932777604035181,0.05556817897035415,0.1278719599607729,1.058063999417129,0.07973214163530429,0.3456803888757634,0.1949898972847421,0.199637989548819,1.979357017700188,3.774483890428684)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowData = $data[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $addr = $cols[$c] + $rowNum
        $ws.Range($addr).Value = $rowData[$c]
    }
}
